# EncuestasAndalucia2 / Hoja1 update:
#  - Insert two new survey rows at the top of the data table (row 2),
#    pushing all existing survey rows down by two.
#  - Row 2: new "Metroscopia" poll (4-15 Nov 2018).
#  - Row 3: new "GAD3" poll (8-14 Nov 2018).
#  - Move the active selection to I2 (VOX column of the newest row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 2 (old row 2 -> row 4, etc.)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# New row 2 - Metroscopia (04/11/2018 - 15/11/2018)
$ws.Range("A2").Value = 43408
$ws.Range("B2").Value = 43419
$ws.Range("C2").Value = 4000
$ws.Range("D2").Value = "Metroscopia"
$ws.Range("E2").Value = 20.9
$ws.Range("F2").Value = 30.9
$ws.Range("G2").Value = 21.2
$ws.Range("H2").Value = 20.1

# New row 3 - GAD3 (08/11/2018 - 14/11/2018)
$ws.Range("A3").Value = 43412
$ws.Range("B3").Value = 43418
$ws.Range("C3").Value = 1803
$ws.Range("D3").Value = "GAD3"
$ws.Range("E3").Value = 17.4
$ws.Range("F3").Value = 33.5
$ws.Range("G3").Value = 19.5
$ws.Range("H3").Value = 22.3
$ws.Range("I3").Value = 3.6

# Match the author's final selection/cursor position.
$ws.Range("I2").Select()
